# regen sval data to filter save games
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @{ B = 3.230985683306322;  C = 1.667794583268128; D = 0.1575252929769615; E = 0.496779210170732;  G = 5.553084769722144 }
    3  = @{ B = 0.6753301551942219; C = 1.667794583268128; D = 0.8054896365839992; E = 0.496779210170732;  G = 3.645393585217082 }
    4  = @{ B = 0.3048080303191223; C = 1.667794583268128; D = 0.1575252929769615; E = 0.496779210170732;  G = 2.626907116734944 }
    5  = @{ B = 3.230985683306322;  C = 1.667794583268128; D = 3.900430680208489;  E = 0.496779210170732;  G = 9.295990156953671 }
    6  = @{ B = 0.6753301551942219; C = 1.667794583268128; D = 26.21740644021617;  E = 8.660232485948974;  G = 37.2207636646275 }
    7  = @{ B = 3.230985683306322;  C = 1.667794583268128; D = 26.21740644021617;  E = 8.660232485948974;  G = 39.7764191927396 }
    8  = @{ B = 3.230985683306322;  C = 1.667794583268128; D = 0.8054896365839992; E = 0.496779210170732;  G = 6.201049113329182 }
    9  = @{ B = 3.230985683306322;  C = 1.667794583268128; D = 3.900430680208489;  E = 8.660232485948974;  G = 17.45944343273191 }
    10 = @{ B = 3.230985683306322;  C = 1.667794583268128; D = 3.900430680208489;  E = 0.496779210170732;  G = 9.295990156953671 }
    11 = @{ B = 0.127881588408715;  C = 0.3127903958511391; D = 0.1575252929769615; E = 0.496779210170732; G = 1.094976487407548 }
    12 = @{ B = 1.459612070389937;  C = 1.667794583268128; D = 0.8054896365839992; E = 0.496779210170732;  G = 4.429675500412797 }
    13 = @{ B = 0.04763786555579896; C = 0.002777888934908601; D = 0.8054896365839992; E = 0.496779210170732; G = 1.352684601245439 }
}

foreach ($r in $data.Keys) {
    $row = $data[$r]
    $ws.Range("B$r").Value = $row.B
    $ws.Range("C$r").Value = $row.C
    $ws.Range("D$r").Value = $row.D
    $ws.Range("E$r").Value = $row.E
    $ws.Range("G$r").Value = $row.G
}
